$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: 2215. Find the Difference of Two Arrays ---
$ws.Cells.Item(28, 1).Value = "2215. Find the Difference of Two Arrays"
$ws.Cells.Item(28, 2).Value = "Easy"
$ws.Cells.Item(28, 2).Interior.Color = 5287936
$ws.Cells.Item(28, 3).Value = "Hashing"
$ws.Cells.Item(28, 4).Value = "Symmetric Difference with HashSets. Create and populate Hashlists, convert to ArrayLists, then use removeAll() method to find the symmetric difference."

$link28 = "https://leetcode.com/problems/find-the-difference-of-two-arrays/solutions/4390106/97-beats-only-using-hashset-friendly/?envType=study-plan-v2&envId=leetcode-75"
$cellE28 = $ws.Cells.Item(28, 5)
$cellE28.Value = $link28 + " "
$ws.Hyperlinks.Add($cellE28, $link28)
$cellE28.Style = "Hyperlink"

# --- Row 29: 2095. Delete the Middle Node of a Linked List ---
$ws.Cells.Item(29, 1).Value = "2095. Delete the Middle Node of a Linked List"
$ws.Cells.Item(29, 2).Value = "Medium"
$ws.Cells.Item(29, 2).Interior.Color = 49407
$ws.Cells.Item(29, 3).Value = "Linked List"
$ws.Cells.Item(29, 4).Value = "1st pass to find the length, calculate the mid point, then traverse and connect prev.next to prev.next.next, or to slow to fast, skipping the nth node."

$link29 = "https://leetcode.com/problems/delete-the-middle-node-of-a-linked-list/solutions/4335889/100-best-approach-slow-fast-pointers/?envType=study-plan-v2&envId=leetcode-75"
$cellE29 = $ws.Cells.Item(29, 5)
$cellE29.Value = $link29 + " "
$ws.Hyperlinks.Add($cellE29, $link29)
$cellE29.Style = "Hyperlink"

# --- View state: selection moves to D32 (and the saved view naturally re-anchors to A1) ---
$ws.Range("D32").Select()
